$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Finished blank rinse samples: update the "Blank" sample's (row 14) result
# from Pass -> Fail, and its A# ug/mL reading from 0.16 -> 25.
$ws.Range("D14").Value = "Fail"
$ws.Range("E14").Value = 25

# Echo the "Other peak(s)" note (same text already present in F14) just
# above the results table as well, matching the header row's formatting.
$ws.Range("A13").Copy()
$ws.Range("F10").PasteSpecial(-4122)
$ws.Range("F10").Value = "6.23, 2.15, 0.015, 8.23"

# Scroll the view down a bit and leave the "Other peak(s)" cell for the
# blank sample selected.
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F14").Select()
